$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.770.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.19%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.649.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.49%  "

# Row 4
$ws.Range("E4").Value = "  +0.45%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.70%  "

# Row 6
$ws.Range("E6").Value = "  +1.61%  "

# Row 7
$ws.Range("E7").Value = "  +0.30%  "

# Row 8
$ws.Range("E8").Value = "  +1.64%  "

# Row 9
$ws.Range("E9").Value = "  +0.81%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.34%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.28%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.878.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.41%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.681.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.21%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.61%  "

# Row 15
$ws.Range("E15").Value = "  +2.06%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.10%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.786.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.11%  "

# Row 18
$ws.Range("E18").Value = "  +0.78%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.43%  "

# Row 20
$ws.Range("E20").Value = "  +0.36%  "

# Row 21
$ws.Range("E21").Value = "  +2.02%  "

# Row 22
$ws.Range("E22").Value = "  +0.63%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +17.20%  "

# Row 24
$ws.Range("E24").Value = "  +2.58%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.54%  "

# Row 26
$ws.Range("E26").Value = "  +0.37%  "

# Row 27
$ws.Range("E27").Value = "  +0.62%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.29%  "

# Row 29
$ws.Range("E29").Value = "  +1.58%  "

# Row 30
$ws.Range("E30").Value = "  +1.73%  "

# Row 31
$ws.Range("E31").Value = "  +1.67%  "

# Row 32
$ws.Range("E32").Value = "  +1.27%  "

# Row 33
$ws.Range("E33").Value = "  +2.65%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.282.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.12%  "

# Row 35
$ws.Range("E35").Value = "  +3.74%  "

# Row 36
$ws.Range("E36").Value = "  +1.99%  "

# Row 37
$ws.Range("E37").Value = "  +3.26%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.536"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.14%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.828"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.30%  "

# Row 40
$ws.Range("E40").Value = "  +0.37%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.815"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.18%  "

# Row 42
$ws.Range("E42").Value = "  -0.47%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.39%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.788.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.49%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.10%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.34%  "

# Row 47
$ws.Range("E47").Value = "  +1.85%  "

# Row 48
$ws.Range("E48").Value = "  +1.30%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.26%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0970"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.14%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.408"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.34%  "
